$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 was previously empty (data started at row 2). Populate it with header labels.
$ws.Range("A1").Value = "Applied Language"
$ws.Range("B1").Value = "Expected Data"
$ws.Range("B1").WrapText = $true

# Fix selection / view to match target state
$ws.Range("B1").Select()
